$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns present in the new rows (template copy skips D to avoid inheriting column D style)
$templateCols = @(1,2,3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)

for ($r = 482; $r -le 495; $r++) {
    foreach ($c in $templateCols) {
        $src = $ws.Cells.Item(2, $c)
        $dst = $ws.Cells.Item($r, $c)
        $src.Copy($dst)
    }
}

# Row 482
$ws.Cells.Item(482,2).Value = 45913
$ws.Cells.Item(482,3).Value = "Global"
$ws.Cells.Item(482,5).Value = "Omar Benyounes"
$ws.Cells.Item(482,6).Value = "center midfield"
$ws.Cells.Item(482,7).Value = "01:14:46"
$ws.Cells.Item(482,8).Value = 9.14
$ws.Cells.Item(482,9).Value = 2.01
$ws.Cells.Item(482,10).Value = 7.1
$ws.Cells.Item(482,11).Value = 1.39
$ws.Cells.Item(482,12).Value = 0.57
$ws.Cells.Item(482,13).Value = 0.07
$ws.Cells.Item(482,14).Value = 0
$ws.Cells.Item(482,15).Value = 7
$ws.Cells.Item(482,16).Value = 7.31
$ws.Cells.Item(482,17).Value = 28.65
$ws.Cells.Item(482,18).Value = 4.75
$ws.Cells.Item(482,19).Value = 34
$ws.Cells.Item(482,20).Value = 6
$ws.Cells.Item(482,21).Value = 40
$ws.Cells.Item(482,22).Value = 8

# Row 483
$ws.Cells.Item(483,2).Value = 45913
$ws.Cells.Item(483,3).Value = "Global"
$ws.Cells.Item(483,5).Value = "Jeremie Laurent"
$ws.Cells.Item(483,6).Value = "left forward"
$ws.Cells.Item(483,7).Value = "01:14:17"
$ws.Cells.Item(483,8).Value = 8.44
$ws.Cells.Item(483,9).Value = 1.84
$ws.Cells.Item(483,10).Value = 6.57
$ws.Cells.Item(483,11).Value = 0.99
$ws.Cells.Item(483,12).Value = 0.62
$ws.Cells.Item(483,13).Value = 0.25
$ws.Cells.Item(483,14).Value = 0.01
$ws.Cells.Item(483,15).Value = 18
$ws.Cells.Item(483,16).Value = 6.8
$ws.Cells.Item(483,17).Value = 30.79
$ws.Cells.Item(483,18).Value = 4.8
$ws.Cells.Item(483,19).Value = 51
$ws.Cells.Item(483,20).Value = 13
$ws.Cells.Item(483,21).Value = 38
$ws.Cells.Item(483,22).Value = 17

# Row 484
$ws.Cells.Item(484,2).Value = 45913
$ws.Cells.Item(484,3).Value = "Global"
$ws.Cells.Item(484,5).Value = "Amir Etien"
$ws.Cells.Item(484,6).Value = "right forward"
$ws.Cells.Item(484,7).Value = "00:19:28"
$ws.Cells.Item(484,8).Value = 1.71
$ws.Cells.Item(484,9).Value = 0.48
$ws.Cells.Item(484,10).Value = 1.23
$ws.Cells.Item(484,11).Value = 0.21
$ws.Cells.Item(484,12).Value = 0.12
$ws.Cells.Item(484,13).Value = 0.15
$ws.Cells.Item(484,14).Value = 0
$ws.Cells.Item(484,15).Value = 9
$ws.Cells.Item(484,16).Value = 5.18
$ws.Cells.Item(484,17).Value = 29.08
$ws.Cells.Item(484,18).Value = 4.48
$ws.Cells.Item(484,19).Value = 5
$ws.Cells.Item(484,20).Value = 5
$ws.Cells.Item(484,21).Value = 4
$ws.Cells.Item(484,22).Value = 3

# Row 485
$ws.Cells.Item(485,2).Value = 45913
$ws.Cells.Item(485,3).Value = "Global"
$ws.Cells.Item(485,5).Value = "Emmanuel Valey"
$ws.Cells.Item(485,6).Value = "left forward"
$ws.Cells.Item(485,7).Value = "00:19:50"
$ws.Cells.Item(485,8).Value = 2.05
$ws.Cells.Item(485,9).Value = 0.5
$ws.Cells.Item(485,10).Value = 1.55
$ws.Cells.Item(485,11).Value = 0.27
$ws.Cells.Item(485,12).Value = 0.12
$ws.Cells.Item(485,13).Value = 0.1
$ws.Cells.Item(485,14).Value = 0
$ws.Cells.Item(485,15).Value = 5
$ws.Cells.Item(485,16).Value = 6.17
$ws.Cells.Item(485,17).Value = 30.24
$ws.Cells.Item(485,18).Value = 4.92
$ws.Cells.Item(485,19).Value = 7
$ws.Cells.Item(485,20).Value = 3
$ws.Cells.Item(485,21).Value = 6
$ws.Cells.Item(485,22).Value = 5

# Row 486
$ws.Cells.Item(486,2).Value = 45913
$ws.Cells.Item(486,3).Value = "Global"
$ws.Cells.Item(486,5).Value = "Sofiane Belle"
$ws.Cells.Item(486,6).Value = "left forward"
$ws.Cells.Item(486,7).Value = "00:33:45"
$ws.Cells.Item(486,8).Value = 2.98
$ws.Cells.Item(486,9).Value = 0.65
$ws.Cells.Item(486,10).Value = 2.32
$ws.Cells.Item(486,11).Value = 0.33
$ws.Cells.Item(486,12).Value = 0.25
$ws.Cells.Item(486,13).Value = 0.08
$ws.Cells.Item(486,14).Value = 0
$ws.Cells.Item(486,15).Value = 7
$ws.Cells.Item(486,16).Value = 4.81
$ws.Cells.Item(486,17).Value = 29.41
$ws.Cells.Item(486,18).Value = 4.68
$ws.Cells.Item(486,19).Value = 5
$ws.Cells.Item(486,20).Value = 2
$ws.Cells.Item(486,21).Value = 14
$ws.Cells.Item(486,22).Value = 2

# Row 487
$ws.Cells.Item(487,2).Value = 45913
$ws.Cells.Item(487,3).Value = "Global"
$ws.Cells.Item(487,5).Value = "Malik Boussaid"
$ws.Cells.Item(487,6).Value = "left back"
$ws.Cells.Item(487,7).Value = "01:21:47"
$ws.Cells.Item(487,8).Value = 8.55
$ws.Cells.Item(487,9).Value = 1.65
$ws.Cells.Item(487,10).Value = 6.87
$ws.Cells.Item(487,11).Value = 0.96
$ws.Cells.Item(487,12).Value = 0.48
$ws.Cells.Item(487,13).Value = 0.22
$ws.Cells.Item(487,14).Value = 0.01
$ws.Cells.Item(487,15).Value = 11
$ws.Cells.Item(487,16).Value = 6.16
$ws.Cells.Item(487,17).Value = 31.06
$ws.Cells.Item(487,18).Value = 4.4
$ws.Cells.Item(487,19).Value = 26
$ws.Cells.Item(487,20).Value = 2
$ws.Cells.Item(487,21).Value = 24
$ws.Cells.Item(487,22).Value = 12
# apply style 6 to E487 (copy format only from an existing styled cell, then restore value)
$styleSrc = $ws.Range("E468")
$eCell = $ws.Cells.Item(487,5)
$styleSrc.Copy($eCell)
$eCell.Value = "Malik Boussaid"

# Row 488
$ws.Cells.Item(488,2).Value = 45913
$ws.Cells.Item(488,3).Value = "Global"
$ws.Cells.Item(488,5).Value = "Ilyes Boughanmi"
$ws.Cells.Item(488,6).Value = "center forward"
$ws.Cells.Item(488,7).Value = "01:34:21"
$ws.Cells.Item(488,8).Value = 8.29
$ws.Cells.Item(488,9).Value = 1.43
$ws.Cells.Item(488,10).Value = 6.84
$ws.Cells.Item(488,11).Value = 0.88
$ws.Cells.Item(488,12).Value = 0.39
$ws.Cells.Item(488,13).Value = 0.15
$ws.Cells.Item(488,14).Value = 0.03
$ws.Cells.Item(488,15).Value = 12
$ws.Cells.Item(488,16).Value = 5.26
$ws.Cells.Item(488,17).Value = 31.5
$ws.Cells.Item(488,18).Value = 4.28
$ws.Cells.Item(488,19).Value = 28
$ws.Cells.Item(488,20).Value = 1
$ws.Cells.Item(488,21).Value = 33
$ws.Cells.Item(488,22).Value = 15

# Row 489
$ws.Cells.Item(489,2).Value = 45913
$ws.Cells.Item(489,3).Value = "Global"
$ws.Cells.Item(489,5).Value = "Yoan Zouma"
$ws.Cells.Item(489,6).Value = "center back"
$ws.Cells.Item(489,7).Value = "01:34:15"
$ws.Cells.Item(489,8).Value = 8.21
$ws.Cells.Item(489,9).Value = 0.88
$ws.Cells.Item(489,10).Value = 7.32
$ws.Cells.Item(489,11).Value = 0.62
$ws.Cells.Item(489,12).Value = 0.21
$ws.Cells.Item(489,13).Value = 0.06
$ws.Cells.Item(489,14).Value = 0
$ws.Cells.Item(489,15).Value = 3
$ws.Cells.Item(489,16).Value = 5.19
$ws.Cells.Item(489,17).Value = 29.78
$ws.Cells.Item(489,18).Value = 4.54
$ws.Cells.Item(489,19).Value = 30
$ws.Cells.Item(489,20).Value = 6
$ws.Cells.Item(489,21).Value = 22
$ws.Cells.Item(489,22).Value = 5

# Row 490
$ws.Cells.Item(490,2).Value = 45913
$ws.Cells.Item(490,3).Value = "Global"
$ws.Cells.Item(490,5).Value = "Karahali Souaré"
$ws.Cells.Item(490,6).Value = "right forward"
$ws.Cells.Item(490,7).Value = "01:02:27"
$ws.Cells.Item(490,8).Value = 6.68
$ws.Cells.Item(490,9).Value = 1.28
$ws.Cells.Item(490,10).Value = 5.39
$ws.Cells.Item(490,11).Value = 0.8
$ws.Cells.Item(490,12).Value = 0.36
$ws.Cells.Item(490,13).Value = 0.13
$ws.Cells.Item(490,14).Value = 0
$ws.Cells.Item(490,15).Value = 11
$ws.Cells.Item(490,16).Value = 6.37
$ws.Cells.Item(490,17).Value = 28.53
$ws.Cells.Item(490,18).Value = 4.8
$ws.Cells.Item(490,19).Value = 45
$ws.Cells.Item(490,20).Value = 9
$ws.Cells.Item(490,21).Value = 22
$ws.Cells.Item(490,22).Value = 9

# Row 491
$ws.Cells.Item(491,2).Value = 45913
$ws.Cells.Item(491,3).Value = "Global"
$ws.Cells.Item(491,5).Value = "Levy Ndoutoume"
$ws.Cells.Item(491,6).Value = "left back"
$ws.Cells.Item(491,7).Value = "01:34:14"
$ws.Cells.Item(491,8).Value = 8.81
$ws.Cells.Item(491,9).Value = 1.33
$ws.Cells.Item(491,10).Value = 7.46
$ws.Cells.Item(491,11).Value = 0.88
$ws.Cells.Item(491,12).Value = 0.37
$ws.Cells.Item(491,13).Value = 0.11
$ws.Cells.Item(491,14).Value = 0
$ws.Cells.Item(491,15).Value = 8
$ws.Cells.Item(491,16).Value = 5.55
$ws.Cells.Item(491,17).Value = 29.19
$ws.Cells.Item(491,18).Value = 4.92
$ws.Cells.Item(491,19).Value = 47
$ws.Cells.Item(491,20).Value = 19
$ws.Cells.Item(491,21).Value = 42
$ws.Cells.Item(491,22).Value = 19

# Row 492
$ws.Cells.Item(492,2).Value = 45913
$ws.Cells.Item(492,3).Value = "Global"
$ws.Cells.Item(492,5).Value = "Naim Ighbane"
$ws.Cells.Item(492,6).Value = "center back"
$ws.Cells.Item(492,7).Value = "01:34:29"
$ws.Cells.Item(492,8).Value = 8.59
$ws.Cells.Item(492,9).Value = 1.02
$ws.Cells.Item(492,10).Value = 7.55
$ws.Cells.Item(492,11).Value = 0.79
$ws.Cells.Item(492,12).Value = 0.2
$ws.Cells.Item(492,13).Value = 0.04
$ws.Cells.Item(492,14).Value = 0
$ws.Cells.Item(492,15).Value = 2
$ws.Cells.Item(492,16).Value = 5.39
$ws.Cells.Item(492,17).Value = 29.07
$ws.Cells.Item(492,18).Value = 4.85
$ws.Cells.Item(492,19).Value = 30
$ws.Cells.Item(492,20).Value = 3
$ws.Cells.Item(492,21).Value = 29
$ws.Cells.Item(492,22).Value = 4

# Row 493
$ws.Cells.Item(493,2).Value = 45913
$ws.Cells.Item(493,3).Value = "Global"
$ws.Cells.Item(493,5).Value = "Hedi Nasri"
$ws.Cells.Item(493,6).Value = "right back"
$ws.Cells.Item(493,7).Value = "00:12:41"
$ws.Cells.Item(493,8).Value = 1.04
$ws.Cells.Item(493,9).Value = 0.07
$ws.Cells.Item(493,10).Value = 0.97
$ws.Cells.Item(493,11).Value = 0.05
$ws.Cells.Item(493,12).Value = 0.02
$ws.Cells.Item(493,13).Value = 0
$ws.Cells.Item(493,14).Value = 0
$ws.Cells.Item(493,15).Value = 0
$ws.Cells.Item(493,16).Value = 4.74
$ws.Cells.Item(493,17).Value = 22.5
$ws.Cells.Item(493,18).Value = 3.89
$ws.Cells.Item(493,19).Value = 1
$ws.Cells.Item(493,20).Value = 0
$ws.Cells.Item(493,21).Value = 4
$ws.Cells.Item(493,22).Value = 0

# Row 494
$ws.Cells.Item(494,2).Value = 45913
$ws.Cells.Item(494,3).Value = "Global"
$ws.Cells.Item(494,5).Value = "Naim Dhib"
$ws.Cells.Item(494,6).Value = "center midfield"
$ws.Cells.Item(494,7).Value = "01:35:42"
$ws.Cells.Item(494,8).Value = 9.26
$ws.Cells.Item(494,9).Value = 1.28
$ws.Cells.Item(494,10).Value = 7.96
$ws.Cells.Item(494,11).Value = 0.87
$ws.Cells.Item(494,12).Value = 0.33
$ws.Cells.Item(494,13).Value = 0.09
$ws.Cells.Item(494,14).Value = 0
$ws.Cells.Item(494,15).Value = 7
$ws.Cells.Item(494,16).Value = 5.73
$ws.Cells.Item(494,17).Value = 27.98
$ws.Cells.Item(494,18).Value = 4.97
$ws.Cells.Item(494,19).Value = 42
$ws.Cells.Item(494,20).Value = 8
$ws.Cells.Item(494,21).Value = 32
$ws.Cells.Item(494,22).Value = 3

# Row 495
$ws.Cells.Item(495,2).Value = 45913
$ws.Cells.Item(495,3).Value = "Global"
$ws.Cells.Item(495,5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(495,6).Value = "center midfield"
$ws.Cells.Item(495,7).Value = "01:33:45"
$ws.Cells.Item(495,8).Value = 10.82
$ws.Cells.Item(495,9).Value = 1.7
$ws.Cells.Item(495,10).Value = 9.1
$ws.Cells.Item(495,11).Value = 1.34
$ws.Cells.Item(495,12).Value = 0.28
$ws.Cells.Item(495,13).Value = 0.1
$ws.Cells.Item(495,14).Value = 0
$ws.Cells.Item(495,15).Value = 6
$ws.Cells.Item(495,16).Value = 6.89
$ws.Cells.Item(495,17).Value = 29.71
$ws.Cells.Item(495,18).Value = 4.6
$ws.Cells.Item(495,19).Value = 40
$ws.Cells.Item(495,20).Value = 10
$ws.Cells.Item(495,21).Value = 31
$ws.Cells.Item(495,22).Value = 7

# Column A (match name) filled in last, matching shared-string append order
$ws.Cells.Item(482,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(483,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(484,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(485,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(486,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(487,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(488,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(489,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(490,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(491,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(492,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(493,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(494,1).Value = "CDF T3 VS Plaine Tonique (R3)"
$ws.Cells.Item(495,1).Value = "CDF T3 VS Plaine Tonique (R3)"

# Update selection to match the final state (A486)
$ws.Range("A486").Select()
